# chore: update Sheets via scheduled runner
# Refresh market-board derived profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across several job sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) to reflect latest pricing data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 284.53
$ws.Range("I15").Value = 284.53
$ws.Range("K15").Value = 853.5899999999999
$ws.Range("M15").Value = -684.5899999999999

# Row 28
$ws.Range("H28").Value = 964.7059
$ws.Range("I28").Value = 1107.9286
$ws.Range("J28").Value = 296.33334
$ws.Range("K28").Value = 1107.9286
$ws.Range("L28").Value = 296.33334
$ws.Range("M28").Value = -622.9286
$ws.Range("N28").Value = -1266.33334

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# Row 132
$ws.Range("H132").Value = 1430178.5
$ws.Range("I132").Value = 1614.1
$ws.Range("J132").Value = 10001565
$ws.Range("K132").Value = 4842.299999999999
$ws.Range("L132").Value = 30004695
$ws.Range("M132").Value = -2312.299999999999
$ws.Range("N132").Value = -30009755

# Row 137
$ws.Range("H137").Value = 1097.2222
$ws.Range("I137").Value = 995.2
$ws.Range("J137").Value = 1224.75
$ws.Range("K137").Value = 2985.6
$ws.Range("L137").Value = 3674.25
$ws.Range("M137").Value = -435.6000000000004
$ws.Range("N137").Value = -8774.25

# Row 138
$ws.Range("H138").Value = 3099.93
$ws.Range("I138").Value = 1057
$ws.Range("J138").Value = 4352.0483
$ws.Range("K138").Value = 3171
$ws.Range("L138").Value = 13056.1449
$ws.Range("M138").Value = 1969
$ws.Range("N138").Value = -23336.1449

# Row 141
$ws.Range("H141").Value = 2045.875
$ws.Range("I141").Value = 2427.8333
$ws.Range("J141").Value = 900
$ws.Range("K141").Value = 7283.499899999999
$ws.Range("L141").Value = 2700
$ws.Range("M141").Value = -2103.499899999999
$ws.Range("N141").Value = -13060

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2263.78
$ws.Range("I32").Value = 2263.78
$ws.Range("K32").Value = 2263.78
$ws.Range("M32").Value = -1976.78

# Row 45
$ws.Range("H45").Value = 1002.4194
$ws.Range("I45").Value = 903.4091
$ws.Range("K45").Value = 903.4091
$ws.Range("M45").Value = -526.4091

# Row 61
$ws.Range("H61").Value = 1986.8572
$ws.Range("I61").Value = 1986.2
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1986.2
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1774.2
$ws.Range("N61").Value = -2424

# Row 74
$ws.Range("H74").Value = 1081.871
$ws.Range("I74").Value = 1084.6
$ws.Range("K74").Value = 1084.6
$ws.Range("M74").Value = -210.5999999999999

# Row 77
$ws.Range("H77").Value = 1081.871
$ws.Range("I77").Value = 1084.6
$ws.Range("K77").Value = 5423
$ws.Range("M77").Value = -1055

# Row 136
$ws.Range("H136").Value = 1986.8572
$ws.Range("I136").Value = 1986.2
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5958.6
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3408.6
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1457.7556
$ws.Range("I134").Value = 1281.3721
$ws.Range("J134").Value = 5250
$ws.Range("K134").Value = 3844.1163
$ws.Range("L134").Value = 15750
$ws.Range("M134").Value = -1309.1163
$ws.Range("N134").Value = -20820

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 967.4
$ws.Range("I22").Value = 1321.2222
$ws.Range("J22").Value = 436.66666
$ws.Range("K22").Value = 1321.2222
$ws.Range("L22").Value = 436.66666
$ws.Range("M22").Value = -971.2221999999999
$ws.Range("N22").Value = -1136.66666

# Row 31
$ws.Range("H31").Value = 38573.785
$ws.Range("I31").Value = 3141.4546
$ws.Range("J31").Value = 168492.33
$ws.Range("K31").Value = 3141.4546
$ws.Range("L31").Value = 168492.33
$ws.Range("M31").Value = -2846.4546
$ws.Range("N31").Value = -169082.33

# Row 34
$ws.Range("H34").Value = 38573.785
$ws.Range("I34").Value = 3141.4546
$ws.Range("J34").Value = 168492.33
$ws.Range("K34").Value = 3141.4546
$ws.Range("L34").Value = 168492.33
$ws.Range("M34").Value = -2939.4546
$ws.Range("N34").Value = -168896.33

# Row 58
$ws.Range("H58").Value = 2558.8
$ws.Range("I58").Value = 892.53656
$ws.Range("J58").Value = 5405.3335
$ws.Range("K58").Value = 892.53656
$ws.Range("L58").Value = 5405.3335
$ws.Range("M58").Value = -689.53656
$ws.Range("N58").Value = -5811.3335

# Row 132
$ws.Range("H132").Value = 1233.6177
$ws.Range("I132").Value = 1000.7273
$ws.Range("J132").Value = 1660.5834
$ws.Range("K132").Value = 3002.1819
$ws.Range("L132").Value = 4981.7502
$ws.Range("M132").Value = -472.1819
$ws.Range("N132").Value = -10041.7502

# Row 136
$ws.Range("H136").Value = 2558.8
$ws.Range("I136").Value = 892.53656
$ws.Range("J136").Value = 5405.3335
$ws.Range("K136").Value = 2677.60968
$ws.Range("L136").Value = 16216.0005
$ws.Range("M136").Value = -127.60968
$ws.Range("N136").Value = -21316.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 3663.4644
$ws.Range("I3").Value = 875.1177
$ws.Range("K3").Value = 2625.3531
$ws.Range("M3").Value = -2513.3531

# Row 131
$ws.Range("H131").Value = 20041788
$ws.Range("I131").Value = 83501830
$ws.Range("J131").Value = 1773.6316
$ws.Range("K131").Value = 250505490
$ws.Range("L131").Value = 5320.8948
$ws.Range("M131").Value = -250500450
$ws.Range("N131").Value = -15400.8948

# Row 137
$ws.Range("H137").Value = 4680.05
$ws.Range("I137").Value = 3266.923
$ws.Range("J137").Value = 7304.4287
$ws.Range("K137").Value = 9800.769
$ws.Range("L137").Value = 21913.2861
$ws.Range("M137").Value = -4700.769
$ws.Range("N137").Value = -32113.2861

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1256.8572
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 699.5
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 699.5
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -1075.5

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 11185.714
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 11660
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 11660
$ws.Range("M54").Value = -9480
$ws.Range("N54").Value = -12700

# Row 107
$ws.Range("H107").Value = 589.7917
$ws.Range("I107").Value = 567.7646999999999
$ws.Range("J107").Value = 643.2857
$ws.Range("K107").Value = 1703.2941
$ws.Range("L107").Value = 1929.8571
$ws.Range("M107").Value = 216.7059000000002
$ws.Range("N107").Value = -5769.8571

# Row 127
$ws.Range("H127").Value = 42597.25
$ws.Range("J127").Value = 42597.25
$ws.Range("L127").Value = 42597.25
$ws.Range("N127").Value = -52517.25

# Row 132
$ws.Range("H132").Value = 1924.4706
$ws.Range("I132").Value = 1419.75
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 4259.25
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -1729.25
$ws.Range("N132").Value = -35060
